# "Matrice di Tracciabilita" - RAD update
# 1) Fix project-name typo: "MeedQueue" -> "MedQueue"
# 2) Add a condensed second traceability table (cols H:O, rows 14-19) that
#    repeats the ID / Use Case / Sequence Diagram / Object Diagram /
#    Statechart / Componenti del Sistema / Moduli di Sistema / Test Case
#    columns of the main table (cols B, I, J, K, L, M, N, O of rows 7-12).
# 3) Leave the active selection on the title block (B2:F4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "MedQueue"

$srcRows = @(7, 8, 9, 10, 11, 12)
$dstRows = @(14, 15, 16, 17, 18, 19)
$copyCols = @("I", "J", "K", "L", "M", "N", "O")

for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $srcRow = $srcRows[$i]
    $dstRow = $dstRows[$i]

    # Column B of the original table becomes column H of the new table.
    $ws.Range("B$srcRow").Copy($ws.Range("H$dstRow"))

    # Columns I-O keep the same column letters, just moved to the new rows.
    foreach ($col in $copyCols) {
        $ws.Range("$col$srcRow").Copy($ws.Range("$col$dstRow"))
    }
}

$ws.Range("B2:F4").Select() | Out-Null
